# ----------------------------------------------------------------------------
# Add Traits sheet, populate Answers sheet Add/Minus Traits columns + trait
# lookup helper column, and add a bold header style for the new sheet.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$answers = $wb.Worksheets.Item("Answers")

# --- Populate "Add Traits" (C) / "Minus Traits" (D) columns on the Answers sheet
$answers.Range("C2").Value = "Docile, Docile, Observant"
$answers.Range("D2").Value = "Imaginative"
$answers.Range("C3").Value = "Hardy, Observant, Observant"
$answers.Range("C4").Value = "Jolly"
$answers.Range("D4").Value = "Observant"
$answers.Range("C5").Value = "Lonely, Imaginative, Imaginative"
$answers.Range("D5").Value = "Jolly"
$answers.Range("C6").Value = "Lonely, Lonely, Imaginative"
$answers.Range("D6").Value = "Jolly, Hasty"
$answers.Range("C7").Value = "Hardy, Confident"
$answers.Range("D7").Value = "Imaginative"
$answers.Range("C8").Value = "Imaginative, Observant"
$answers.Range("C9").Value = "Imaginative,  Lonely"
$answers.Range("C10").Value = "Imaginative, Imaginative, Jolly"
$answers.Range("D10").Value = "Aggressive"
$answers.Range("C11").Value = "Aggressive, Hardy, Hardy"
$answers.Range("D11").Value = "Confident, Hasty"
$answers.Range("C12").Value = "Confident, Hardy"
$answers.Range("C13").Value = "Confident, Confident"
$answers.Range("D13").Value = "Hardy"
$answers.Range("C14").Value = "Hardy, Hardy, Observant, Aggressive"
$answers.Range("D14").Value = "Lonely"
$answers.Range("C15").Value = "Lonely, Docile, Docile"
$answers.Range("D15").Value = "Hardy, Hardy"
$answers.Range("C16").Value = "Hardy, Hardy, Confident"
$answers.Range("D16").Value = "Hasty, Docile"
$answers.Range("C17").Value = "Docile, Observant, Observant"
$answers.Range("D17").Value = "Jolly"
$answers.Range("C18").Value = "Hasty, Imaginative, Observant"
$answers.Range("D18").Value = "Hardy, Confident"
$answers.Range("C19").Value = "Hasty, Docile, Jolly"
$answers.Range("D19").Value = "Observant"
$answers.Range("C20").Value = "Observant, Confident, Confident, Imaginative"
$answers.Range("D20").Value = "Hasty "
$answers.Range("C21").Value = "Hasty, Docile, Docile, "
$answers.Range("D21").Value = "Imaginative, Jolly"
$answers.Range("C22").Value = "Lonely, Imaginative,  Jolly"
$answers.Range("D22").Value = "Hardy, Aggressive"
$answers.Range("C23").Value = "Hardy, Observant"
$answers.Range("D23").Value = "Jolly, Hasty"
$answers.Range("C24").Value = "Aggressive, Hardy "
$answers.Range("C25").Value = "Docile, Docile"
$answers.Range("D25").Value = "Jolly"
$answers.Range("C26").Value = "Jolly, Jolly, Confident"
$answers.Range("D26").Value = "Observant"
$answers.Range("C27").Value = "Observant, Hardy"
$answers.Range("D27").Value = "Confident  "
$answers.Range("C28").Value = "Confident, Confident, Confident, Hardy"
$answers.Range("D28").Value = "Observant"
$answers.Range("C29").Value = "Docile, Docile, Imaginative, Observant"
$answers.Range("C30").Value = "Confident, Confident, Hardy, Hardy"
$answers.Range("D30").Value = "Observant"
$answers.Range("C31").Value = "Imaginative, Jolly, Hasty"
$answers.Range("D31").Value = "Confident"
$answers.Range("C32").Value = "Docile, Confident"
$answers.Range("D32").Value = "Imagination"
$answers.Range("C33").Value = "Imaginative, Aggressive"
$answers.Range("D33").Value = "Docile"
$answers.Range("C34").Value = "Lonely, Lonely, Imaginative"
$answers.Range("D34").Value = "Docile"

# --- Small trait-name helper/reference list in column H
$answers.Range("H3").Value = "Observant"
$answers.Range("H4").Value = "Confident"
$answers.Range("H5").Value = "Jolly"
$answers.Range("H6").Value = "Hardy"
$answers.Range("H7").Value = "Aggressive"
$answers.Range("H8").Value = "Hasty"
$answers.Range("H9").Value = "Lonely"
$answers.Range("H10").Value = "Docile"
$answers.Range("H11").Value = "Imaginative"

# --- Resize columns C and H on the Answers sheet to fit the new content
$answers.Columns.Item(3).ColumnWidth = 42.57
$answers.Columns.Item(8).ColumnWidth = 16.3

# --- Add the new "Traits" worksheet after "Answers" (end of the workbook)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$traits = $wb.Worksheets.Add($null, $lastSheet)
$traits.Name = "Traits"

# --- Header row (bold)
$traits.Range("B3").Value = "Trait Name"
$traits.Range("C3").Value = "Associated Film Genre"
$traits.Range("D3").Value = "Associated Game Genre"
$traits.Range("B3:D3").Font.Bold = $true

# --- Trait / genre data rows
$traits.Range("B4").Value = "Observant"
$traits.Range("C4").Value = "Mystery, Crime, Sci-fi"
$traits.Range("D4").Value = "Strategy, Role-playing"
$traits.Range("B5").Value = "Confident"
$traits.Range("C5").Value = "Superhero, Adventure, Action-Adventure, Horror"
$traits.Range("D5").Value = "Sports, Adventure"
$traits.Range("B6").Value = "Jolly"
$traits.Range("C6").Value = "Comedy, Animation, Fantasy"
$traits.Range("D6").Value = "Adventure"
$traits.Range("B7").Value = "Hardy"
$traits.Range("C7").Value = "Thriller, Crime"
$traits.Range("D7").Value = "Action, Strategy"
$traits.Range("B8").Value = "Aggressive"
$traits.Range("C8").Value = "Drama, Action, Thriller"
$traits.Range("D8").Value = "Action, Shooter"
$traits.Range("B9").Value = "Hasty"
$traits.Range("C9").Value = "Superhero, Action-comedy, Horror"
$traits.Range("D9").Value = "Action-Adventure, Horror"
$traits.Range("B10").Value = "Lonely"
$traits.Range("C10").Value = "Romance, Comedy, Fantasy"
$traits.Range("D10").Value = "Adventure, Action-Adventure"
$traits.Range("B11").Value = "Docile"
$traits.Range("C11").Value = "Romance"
$traits.Range("D11").Value = "Role-playing, Action"
$traits.Range("B12").Value = "Imaginative"
$traits.Range("C12").Value = "Fantasy, Sci-fi, Animation"
$traits.Range("D12").Value = "Role-playing, Adventure, Action-adventure"

# --- Column widths for the Traits sheet
$traits.Columns.Item(2).ColumnWidth = 22.14
$traits.Columns.Item(3).ColumnWidth = 45.14
$traits.Columns.Item(4).ColumnWidth = 39.29

# --- Selection / view state to mirror the saved workbook
$traits.Range("B4:B12").Select()
$answers.Select()
$answers.Range("E35").Select()

# --- Leave focus on the Answers sheet (matches activeTab of the saved file)
$answers.Activate()
